# Add season-record columns (Wins / Losses / Ties) to the right of the
# existing table, mirroring the author's "get season record" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels in AD1:AF1 -----------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the look of the other header cells (bold, bordered, centered) by
# copying the formatting from the neighboring "Unnamed: 28" header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows (2-56): same record repeated for every player --------------
$ws.Range("AD2:AD56").Value = 80
$ws.Range("AE2:AE56").Value = 82
$ws.Range("AF2:AF56").Value = 0

Write-Output "Added Wins/Losses/Ties columns (AD:AF) for rows 1-56"
